$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper behaviour: Price values that look like plain decimal numbers (single
# "." and only digits) would otherwise be auto-converted to a number by Excel,
# so we force Text format for the assignment and then restore the default
# "Normal" cell style so no stray formatting is left behind.

$ws.Range("D2").Value = '41.526.85'
$ws.Range("E2").Value = '  +0.87%  '

$ws.Range("D3").Value = '2.483.25'
$ws.Range("E3").Value = '  +1.17%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.998'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.18%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '314.49'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.75%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '93.55'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.03%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.543'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -1.33%  '

$ws.Range("E8").Value = '  -0.17%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.509'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +3.71%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '32.79'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.72%  '

$ws.Range("E11").Value = '  +1.74%  '

$ws.Range("E12").Value = '  +2.99%  '

$ws.Range("D13").Value = '2.865.09'
$ws.Range("E13").Value = '  +1.01%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.85'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.90%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '16.07'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +10.51%  '

$ws.Range("D16").Value = '2.519.68'
$ws.Range("E16").Value = '  +1.98%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.763'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.46%  '

$ws.Range("D18").Value = '41.526.51'
$ws.Range("E18").Value = '  +0.86%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.40'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.48%  '

$ws.Range("D20").Value = '0.0₃0937'
$ws.Range("E20").Value = '  +2.70%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '72.04'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +6.84%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '11.44'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +2.58%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '237.51'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.25%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.73'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.22%  '

$ws.Range("E25").Value = '  -0.32%  '

$ws.Range("E26").Value = '  +1.34%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '24.97'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +5.26%  '

$ws.Range("E28").Value = '  -0.27%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.69'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.86%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '36.18'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +2.52%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '157.79'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +4.63%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.50'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.16%  '

$ws.Range("E33").Value = '  +0.53%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0755'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +3.39%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.46'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -7.54%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '17.44'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +4.98%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.94'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.70%  '

$ws.Range("E40").Value = '  +0.85%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '4.12'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.07%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '19.77'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.86%  '

$ws.Range("D44").Value = '1.980.14'
$ws.Range("E44").Value = '  +0.36%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0285'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.31%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.97'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.37%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.09'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +6.93%  '

$ws.Range("D48").Value = '2.722.87'
$ws.Range("E48").Value = '  +1.07%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '98.19'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +2.84%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '68.26'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.87%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '72.57'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.72%  '

# Row 38 and 39 swap: ARBITRUM moves to row 38, Kaspa moves to row 39
$ws.Range("B38").Value = "ARBITRUM"
$ws.Range("C38").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.84'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.61%  '

$ws.Range("B39").Value = "Kaspa"
$ws.Range("C39").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.105'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +2.97%  '